$wb = $excel.ActiveWorkbook

$phase1 = $wb.Worksheets.Item("Phase1")
$phase2 = $wb.Worksheets.Item("Phase2")

$phase1.Range("B2").Value = 0.1
$phase1.Range("B3").Value = 0.1
$phase1.Range("B4").Value = 0.05
$phase1.Range("B5").Value = 0.05

$phase2.Range("C2").Value = 0.1
$phase2.Range("C3").Value = 0.1
$phase2.Range("C4").Value = 0.05
$phase2.Range("C5").Value = 0.05

$phase1.Activate()
$phase1.Range("B2:B5").Select()

$phase2.Activate()
$phase2.Range("C2:C5").Select()
